$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 182 (pushes old rows 182-185 down to 186-189)
$ws.Rows("182:185").Insert()

# Fill in E181 (was blank) and the shared formulas F181:G181 which now extend through row 184
$ws.Range("E181").Value = 0.76041666666666663

# New data rows 182-184 (full entries)
$ws.Range("A182").Value = 2014
$ws.Range("B182").Value = 8
$ws.Range("C182").Value = 12
$ws.Range("D182").Value = 0.4861111111111111
$ws.Range("E182").Value = 0.5

$ws.Range("A183").Value = 2014
$ws.Range("B183").Value = 8
$ws.Range("C183").Value = 12
$ws.Range("D183").Value = 0.64583333333333337
$ws.Range("E183").Value = 0.69791666666666663

$ws.Range("A184").Value = 2014
$ws.Range("B184").Value = 8
$ws.Range("C184").Value = 13
$ws.Range("D184").Value = 0.4201388888888889
$ws.Range("E184").Value = 0.54513888888888895

# New partial row 185 (only A:D filled, E/F left blank, same as old placeholder row pattern)
$ws.Range("A185").Value = 2014
$ws.Range("B185").Value = 8
$ws.Range("C185").Value = 14
$ws.Range("D185").Value = 0.3888888888888889

# Extend shared formulas F/G across the new rows (176:184, matching the diff's shared formula ref)
$ws.Range("F181:F184").Formula = "=(E181-D181)*24*60"
$ws.Range("G181:G184").Formula = "=F181/60"

# Apply number-format styles consistent with the rest of the table to the new C column entries
$ws.Range("C182:C185").NumberFormat = "0"

# Update selection to match the recorded cursor position after the edit
$ws.Range("C183").Select()

$wb.Save()
